$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.620328187942505
$ws.Range("B1").Value = 3.622976303100586
$ws.Range("C1").Value = 2.733088254928589
$ws.Range("D1").Value = 2.499660968780518
$ws.Range("E1").Value = 2.552166700363159
